$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows before row 179, shifting the existing rows
# (old 179..260) down to (181..262).
$ws.Range("A179:A180").EntireRow.Insert()

# New row 179
$ws.Cells.Item(179, 1).Value = 10
$ws.Cells.Item(179, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(179, 3).Value = "La Araucanía"
$ws.Cells.Item(179, 4).Value = 44529
$ws.Cells.Item(179, 5).Value = 9
$ws.Cells.Item(179, 6).Value = 100112037
$ws.Cells.Item(179, 7).Value = "Cebollín"
$ws.Cells.Item(179, 8).Value = "Sin especificar"
$ws.Cells.Item(179, 9).Value = "Primera"
$ws.Cells.Item(179, 10).Value = 125
$ws.Cells.Item(179, 11).Value = 9000
$ws.Cells.Item(179, 12).Value = 9000
$ws.Cells.Item(179, 13).Value = 9000
$ws.Cells.Item(179, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(179, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(179, 16).Value = 750
$ws.Cells.Item(179, 17).Value = 12
$ws.Cells.Item(179, 18).Value = "Hortaliza"

# New row 180
$ws.Cells.Item(180, 1).Value = 10
$ws.Cells.Item(180, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(180, 3).Value = "La Araucanía"
$ws.Cells.Item(180, 4).Value = 44529
$ws.Cells.Item(180, 5).Value = 9
$ws.Cells.Item(180, 6).Value = 100112037
$ws.Cells.Item(180, 7).Value = "Cebollín"
$ws.Cells.Item(180, 8).Value = "Sin especificar"
$ws.Cells.Item(180, 9).Value = "Primera"
$ws.Cells.Item(180, 10).Value = 95
$ws.Cells.Item(180, 11).Value = 5000
$ws.Cells.Item(180, 12).Value = 5000
$ws.Cells.Item(180, 13).Value = 5000
$ws.Cells.Item(180, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(180, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(180, 16).Value = 417
$ws.Cells.Item(180, 17).Value = 12
$ws.Cells.Item(180, 18).Value = "Hortaliza"

Write-Output "done"
